$wb = $excel.ActiveWorkbook

# --- Add the new weekly sheet "2026-02-25" at the end of the tab order ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2026-02-25"

# --- Header row ---
$ws.Cells.Item(1,1).Value = "rank"
$ws.Cells.Item(1,2).Value = "title"
$ws.Cells.Item(1,3).Value = "volume"
$ws.Cells.Item(1,4).Value = "publisher"

# --- Reference cell that already carries the "single/partial-volume" highlight
#     style (light-yellow fill) used throughout the workbook, so we copy that
#     cell's format (reusing its existing style index) onto the new sheet's
#     cells instead of fabricating a brand-new style/fill entry. ---
$styleSheet = $wb.Worksheets.Item("2026-02-18")
$styleSrc = $styleSheet.Cells.Item(15,3)

# --- Data rows: rank <TAB> title <TAB> volume <TAB> highlighted(1/0) ---
$rowsText = @'
1	転生したらスライムだった件	31	0
2	ミステリと言う勿れ	16	0
3	転生賢者の異世界ライフ~第二の職業を得て、世界最強になりました~	31	0
4	チェンソーマン	23	0
5	金色のガッシュ!! 2 Page	37	0
6	転生したら第七王子だったので、気ままに魔術を極めます	22	0
7	信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐&『ざまぁ!』します!	21	0
8	100万の命の上に俺は立っている	23	0
9	「くじ」から始まる婚約生活~厳正なる抽選の結果、笑わない次期公爵様の婚約者に当選しました~	7	0
10	ガチアクタ	18	0
11	BORUTO-ボルト- -TWO BLUE VORTEX-	7	0
12	星降る王国のニナ	18	0
13	最凶貴族は死亡フラグを覆す1	1	1
14	転生したらスライムだった件 異聞 ~魔国暮らしのトリニティ~	13	0
15	出来損ない皇子の成り上がり~聖痕無しの第三皇子に転生したけど、今度こそ家族を守るために最強を目指す~1	1	1
16	身代わりの生贄だったはずの私、凶犬王子の愛に困惑中 1(アリアンローズコミックス)	1	1
17	無能の中の無能王子 スキルを授かりましたが、周りの女性はとかです コミック版	4	0
18	クラス転移に巻き込まれたコンビニ店員のおっさん、勇者には必要なかった余り物スキルを駆使して最強となるようです。 コミック版	1	1
19	辺境の薬師、都でSランク冒険者となる~英雄村の少年がチート薬で無自覚無双~	11	0
20	王様ランキング	21	0
21	水属性の魔法使い@COMIC	8	0
22	マチ姉さんのポンコツおとぎ話アワー	1	1
23	人の恋路を邪魔する僕は	1	1
24	無能の中の無能王子 スキルを授かりましたが、周りの女性はとかです コミック版	3	1
25	金田一パパの事件簿	3	1
26	メダリスト	14	0
27	聖者無双	16	0
28	マチネとソワレ	18	0
29	正直不動産	23	0
30	29歳独身中堅冒険者の日常	21	0
31	桃源暗鬼	28	0
32	BORUTO-ボルト- -TWO BLUE VORTEX-	6	0
33	アオのハコ	24	0
34	大嫌いな運命の人1	1	1
35	悪意がえし	1	1
36	ある日、惰眠を貪っていたら一族から追放されて森に捨てられました そのまま寝てたら周りが勝手に魔物の国を作ってたけど、私は気にせず今日も眠ります コミック版	1	1
37	最下層の僕が奴隷を飼ったら ―監禁観察日記―	1	1
38	転生したら平民でした。~生活水準に耐えられないので貴族を目指します~(コミック)	7	0
39	Sランクパーティから解雇された~『呪いのアイテム』しか作れませんが、その性能はアーティファクト級なり……!~	13	0
40	だれでも抱けるキミが好き	8	0
41	血を這う亡国の王女	6	0
42	鬼の花嫁	8	0
43	ブルーロック	37	0
44	むせるくらいの愛をあげる	7	0
45	おしえて執事くん	2	1
46	ドラハチ	9	0
47	刃牙らへん	6	0
48	BORUTO-ボルト- -TWO BLUE VORTEX-	4	0
49	BORUTO-ボルト- -TWO BLUE VORTEX-	5	0
50	無能令嬢の契約結婚	2	1
51	転生したらスライムだった件 クレイマンREVENGE	8	0
52	英雄と賢者の転生婚~かつての好敵手と婚約して最強夫婦になりました~	7	0
53	悠久の愚者アズリーの、賢者のすゝめ と、ポチの大冒険13	13	0
54	戦隊大失格	21	0
55	蒼く染めろ	19	0
56	魔入りました!入間くん	47	0
57	BORUTO-ボルト- -TWO BLUE VORTEX-	3	1
58	BORUTO-ボルト- -TWO BLUE VORTEX-	1	1
59	BORUTO-ボルト- -TWO BLUE VORTEX-	2	1
60	恋に溺れる配信者	1	1
61	燁姫	1	1
62	無能の中の無能王子 スキルを授かりましたが、周りの女性はとかです コミック版	2	1
63	テイマーさんのVRMMO育成日誌 コミック版	1	1
64	細菌少女	1	1
65	謎のユリイカ	1	1
66	アラフォー賢者の異世界生活日記 ZERO ソード・アンド・ソーサリス・ワールド 第1話	1	1
67	侯爵令嬢は手駒を演じる 8(アリアンローズコミックス)	8	0
68	世界最強の魔女、始めました ~私だけ『攻略サイト』を見れる世界で自由に生きます~	11	0
69	パリピ孔明	24	0
70	幼女とスコップと魔眼王	5	0
71	涙雨とセレナーデ	14	0
72	大自然の魔法師アシュト、廃れた領地でスローライフ7	7	0
73	スーパーの裏でヤニ吸うふたり 通常版	8	0
74	器用貧乏、城を建てる~開拓学園の劣等生なのに、上級職のスキルと魔法がすべて使えます~@COMIC	7	0
75	モンスターがあふれる世界になったので、好きに生きたいと思います	14	0
76	メダリスト	6	0
77	異世界ありがとう	9	0
78	となりの席のヤツがそういう目で見てくる	5	0
79	獣王と薬草	8	0
80	ほどなく、お別れです	5	0
81	メダリスト	7	0
82	メダリスト	8	0
83	メダリスト	9	0
84	メダリスト	10	0
85	生まれた直後に捨てられたけど、前世が大賢者だったので余裕で生きてます ~最強赤ちゃん大暴走~14	14	0
86	恋せよまやかし天使ども	6	0
87	はじめの一歩	145	0
88	青のミブロー新選組編ー	9	0
89	灰仭巫覡	7	0
90	とある魔術の禁書目録	33	0
91	僕の心のヤバイやつ	13	0
92	葬送のフリーレン	15	0
93	BORUTO-ボルト- -NARUTO NEXT GENERATIONS-	20	0
94	BORUTO-ボルト- -NARUTO NEXT GENERATIONS-	19	0
95	BORUTO-ボルト- -NARUTO NEXT GENERATIONS-	17	0
96	BORUTO-ボルト- -NARUTO NEXT GENERATIONS-	18	0
97	BORUTO-ボルト- -NARUTO NEXT GENERATIONS-	5	0
98	BORUTO-ボルト- -NARUTO NEXT GENERATIONS-	6	0
99	BORUTO-ボルト- -NARUTO NEXT GENERATIONS-	7	0
100	逆転エンゲージメント~悪名高い御曹司が私にだけ甘すぎる~	8	0
'@

$lines = $rowsText -split "`r?`n" | Where-Object { $_.Length -gt 0 }

$r = 2
foreach ($line in $lines) {
    $parts = $line -split "`t"
    $rank = [int]$parts[0]
    $title = $parts[1]
    $volume = [int]$parts[2]
    $highlighted = $parts[3] -eq "1"

    $ws.Cells.Item($r,1).Value = $rank
    $ws.Cells.Item($r,2).Value = $title
    $ws.Cells.Item($r,3).Value = $volume

    if ($highlighted) {
        $styleSrc.Copy()
        $ws.Cells.Item($r,3).PasteSpecial(-4122)
    }

    $r++
}

Write-Output "Added sheet '$($ws.Name)' with $($r - 2) data rows; total sheets: $($wb.Worksheets.Count)"
